$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Haras El Hodood - Smouha) odds updates
$ws.Range("S3").Value = 1.57
$ws.Range("T3").Value = 2.3

# Row 4 (Haka - SJK) odds updates
$ws.Range("G4").Value = 3.4
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 3.75
$ws.Range("L4").Value = 2.63
$ws.Range("Q4").Value = 1.67
$ws.Range("R4").Value = 2.15
$ws.Range("U4").Value = 1.57
$ws.Range("V4").Value = 2.25
$ws.Range("X4").Value = 19
$ws.Range("Y4").Value = 12
$ws.Range("AG4").Value = 126
$ws.Range("AH4").Value = 9.5
$ws.Range("AJ4").Value = 9
$ws.Range("AK4").Value = 19
$ws.Range("AO4").Value = 17
$ws.Range("AS4").Value = 126
$ws.Range("AX4").Value = 11

# Row 5 (Skalica - Slovan Bratislava) odds updates
$ws.Range("G5").Value = 6.5
$ws.Range("H5").Value = 4.75
$ws.Range("I5").Value = 1.42
$ws.Range("J5").Value = 5.7
$ws.Range("K5").Value = 2.5
$ws.Range("L5").Value = 1.87
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 1.17
$ws.Range("P5").Value = 4.55
$ws.Range("Q5").Value = 1.55
$ws.Range("R5").Value = 2.42
$ws.Range("S5").Value = 1.29
$ws.Range("T5").Value = 3.4
$ws.Range("U5").Value = 1.72
$ws.Range("V5").Value = 2
$ws.Range("X5").Value = 50
$ws.Range("Y5").Value = 21
$ws.Range("AA5").Value = 65
$ws.Range("AB5").Value = 55
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 10
$ws.Range("AE5").Value = 18
$ws.Range("AF5").Value = 70
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 7.9
$ws.Range("AI5").Value = 8
$ws.Range("AK5").Value = 10.5
$ws.Range("AL5").Value = 11.5
$ws.Range("AM5").Value = 24
$ws.Range("AN5").Value = 7.9
$ws.Range("AP5").Value = 32
$ws.Range("AQ5").Value = 200
$ws.Range("AR5").Value = 200
$ws.Range("AS5").Value = 300
$ws.Range("AT5").Value = 3.4
$ws.Range("AU5").Value = 7.7
$ws.Range("AV5").Value = 60
$ws.Range("AW5").Value = 3.45
$ws.Range("AX5").Value = 6.4
$ws.Range("AZ5").Value = 17
$ws.Range("BB5").Value = 150
